$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "big and tall running pants for men"
$ws.Cells.Item(2, 1).Value = "girls compression knee sleeve"
$ws.Cells.Item(3, 1).Value = "football pants with pads"
$ws.Cells.Item(4, 1).Value = "compression padded shorts"
$ws.Cells.Item(5, 1).Value = "spandex for men pants"
$ws.Cells.Item(6, 1).Value = "elastic knee pad"
$ws.Cells.Item(7, 1).Value = "knee pads gel construction"
$ws.Cells.Item(8, 1).Value = "capri shorts for men"
$ws.Cells.Item(9, 1).Value = "football padded pants"
$ws.Cells.Item(10, 1).Value = "girls volleyball kneepads"
$ws.Cells.Item(11, 1).Value = "black capris men"
$ws.Cells.Item(12, 1).Value = "womens basketball pants"
$ws.Cells.Item(13, 1).Value = "baseball material"
$ws.Cells.Item(14, 1).Value = "youth sports leggings"
$ws.Cells.Item(15, 1).Value = "protect knee pads"
$ws.Cells.Item(16, 1).Value = "knees pad"
$ws.Cells.Item(17, 1).Value = "black youth knee pads"
$ws.Cells.Item(18, 1).Value = "exercise pads for knees"
$ws.Cells.Item(19, 1).Value = "gym shorts above knee for men"
$ws.Cells.Item(20, 1).Value = "knee sleeve basketball youth"
$ws.Cells.Item(21, 1).Value = "knee pads construction"
$ws.Cells.Item(22, 1).Value = "work knee pad"
$ws.Cells.Item(23, 1).Value = "bjj knee"
$ws.Cells.Item(24, 1).Value = "knee sleeve bjj"
$ws.Cells.Item(25, 1).Value = "knee pads under"
$ws.Cells.Item(26, 1).Value = "mens long cycling pants"
$ws.Cells.Item(27, 1).Value = "mens pad"
$ws.Cells.Item(28, 1).Value = "boys long baseball pants"
$ws.Cells.Item(29, 1).Value = "mens basketball gear"
$ws.Cells.Item(30, 1).Value = "girl sliding shorts"
$ws.Cells.Item(31, 1).Value = "calf silicone pads"
$ws.Cells.Item(32, 1).Value = "compression knee sleeve men pair"
$ws.Cells.Item(33, 1).Value = "girls black softball pants"
$ws.Cells.Item(34, 1).Value = "hip protector pads"
$ws.Cells.Item(35, 1).Value = "flexible work pants for men"
$ws.Cells.Item(36, 1).Value = "knee protection pads"
$ws.Cells.Item(37, 1).Value = "lightweight work pants for men"
$ws.Cells.Item(38, 1).Value = "youth baseball pants long"
$ws.Cells.Item(39, 1).Value = "knees pads work"
$ws.Cells.Item(40, 1).Value = "tights boys"
$ws.Cells.Item(41, 1).Value = "mens shorts long below knee"
$ws.Cells.Item(42, 1).Value = "knee sleeve youth"
$ws.Cells.Item(43, 1).Value = "snowboarding pants men"
$ws.Cells.Item(44, 1).Value = "baseball shorts for men"
$ws.Cells.Item(45, 1).Value = "baseball mens pants"
$ws.Cells.Item(46, 1).Value = "knee compression sleeve - reduce strain & swelling"
$ws.Cells.Item(47, 1).Value = "pads men"
$ws.Cells.Item(48, 1).Value = "basketball sleeve youth leg"
$ws.Cells.Item(49, 1).Value = "thigh pads football"
$ws.Cells.Item(50, 1).Value = "compression volleyball"
$ws.Cells.Item(51, 1).Value = "leggings for mens"
$ws.Cells.Item(52, 1).Value = "mens yoga pants"
$ws.Cells.Item(53, 1).Value = "padded football pants"
$ws.Cells.Item(54, 1).Value = "spandex capris"
$ws.Cells.Item(55, 1).Value = "water knee hockey"
$ws.Cells.Item(56, 1).Value = "compression pants sleeves"
$ws.Cells.Item(57, 1).Value = "knee sleeve padded"
$ws.Cells.Item(58, 1).Value = "knees pads for construction"
$ws.Cells.Item(59, 1).Value = "tight capri"
$ws.Cells.Item(60, 1).Value = "mens baseball compression shorts"
$ws.Cells.Item(61, 1).Value = "mens running knee compression"
$ws.Cells.Item(62, 1).Value = "black football leggings"
$ws.Cells.Item(63, 1).Value = "knee sleeves basketball youth"
$ws.Cells.Item(64, 1).Value = "cycling knee pads"
$ws.Cells.Item(65, 1).Value = "construction knee pad"
$ws.Cells.Item(66, 1).Value = "compression calf leggings"
$ws.Cells.Item(67, 1).Value = "baseball youth compression sleeve"
$ws.Cells.Item(68, 1).Value = "6 pairs of leggings"
$ws.Cells.Item(69, 1).Value = "basketball shorts for men pack of 5"
$ws.Cells.Item(70, 1).Value = "compression pants youth boys"
$ws.Cells.Item(71, 1).Value = "impact shorts men"
$ws.Cells.Item(72, 1).Value = "large knee pad"
$ws.Cells.Item(73, 1).Value = "baseball compression sleeve"
$ws.Cells.Item(74, 1).Value = "boys sports leggings"
$ws.Cells.Item(75, 1).Value = "volleyball spandex pack"
$ws.Cells.Item(76, 1).Value = "baseball pants youth large"
$ws.Cells.Item(77, 1).Value = "boys paintball pants"
$ws.Cells.Item(78, 1).Value = "yoga pants mens"
$ws.Cells.Item(79, 1).Value = "calf tear compression sleeve"
$ws.Cells.Item(80, 1).Value = "compression shorts men long length"
$ws.Cells.Item(81, 1).Value = "yoga hand pads"
$ws.Cells.Item(82, 1).Value = "knee sleeves with padding"
$ws.Cells.Item(83, 1).Value = "athletic capri leggings"
$ws.Cells.Item(84, 1).Value = "pants compression men"
$ws.Cells.Item(85, 1).Value = "basketball padding"
$ws.Cells.Item(86, 1).Value = "knee pads for men floor work"
$ws.Cells.Item(87, 1).Value = "youth knee sleeve wrestling"
$ws.Cells.Item(88, 1).Value = "professional construction knee pads"
$ws.Cells.Item(89, 1).Value = "youth basketball"
$ws.Cells.Item(90, 1).Value = "basketball compression knee sleeve"
$ws.Cells.Item(91, 1).Value = "black softball pants youth girls"
$ws.Cells.Item(92, 1).Value = "hex gear wash"
$ws.Cells.Item(93, 1).Value = "knee construction pads"
$ws.Cells.Item(94, 1).Value = "youth girls softball pants"
$ws.Cells.Item(95, 1).Value = "compression tight pants"
$ws.Cells.Item(96, 1).Value = "male workout leggings"
$ws.Cells.Item(97, 1).Value = "boys xl baseball pants"
$ws.Cells.Item(98, 1).Value = "thick leggings for men"
$ws.Cells.Item(99, 1).Value = "knee pads for"
$ws.Cells.Item(100, 1).Value = "adult football girdle"
